$d = $word.ActiveDocument

# The document ends with a list paragraph:
#   "Переделать USB: убрать резисторы, добавить провод 5В."
# followed by the (hidden) _GoBack bookmark. We need to insert two new
# list-item paragraphs after it, and keep the _GoBack bookmark attached
# at the very end of the document (after the newly-added text).

$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$text1 = "Пересмотреть питание и зарядку. А то зарядка может не идти, если прочее будет потреблять много."
$text2 = "Переделать инфракрасный светодиод."

# Use InsertAfter (not InsertBefore) so the new runs correctly inherit the
# surrounding ru-RU character formatting.
$insertionPoint.InsertAfter("`r" + $text1 + "`r" + $text2)

# InsertAfter leaves the pre-existing _GoBack bookmark anchored to its old
# spot (right after "5В."), but it needs to move to the very end of the
# document, after the text we just added.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create the bookmark at the end of the document. Directly targeting a
# zero-length range sitting right at "end-of-paragraph-content" is
# unreliable, so append a throwaway marker character first, bookmark just
# before it (a stable position), then remove the marker again - the
# bookmark stays put.
$finalPara = $d.Paragraphs.Last
$markerPos = $finalPara.Range.End - 1
$markerRange = $d.Range($markerPos, $markerPos)
$markerRange.InsertAfter("X")

$bookmarkRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$markerRange2 = $d.Range($markerPos, $markerPos + 1)
$markerRange2.Delete()
